$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1878.6154
$ws.Cells.Item(19, 8).Value = 9999
